$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before AT, shifting AT:AU -> AU:AV
$ws.Columns("AT:AT").Insert()
$ws.Columns("AT:AT").ColumnWidth = 38.33

# Rename the "Hang khuyen mai" header to "Tinh chat hang hoa"
$ws.Range("AE4").Value = "Tính chất hàng hóa"

# Give the newly inserted column its header text
$ws.Range("AT4").Value = "Thông tin hóa đơn liên quan"

# Update the active selection (also drops the stale scroll position)
$ws.Range("AT17").Select() | Out-Null
